$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1.62
$ws.Range("I4").Value = 6.25
$ws.Range("L4").Value = 6.5
$ws.Range("X4").Value = 6.5
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 21
$ws.Range("AO4").Value = 8.5
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201
$ws.Range("K5").Value = 2.12
$ws.Range("L5").Value = 7.8
$ws.Range("N5").Value = 6.8
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 2.85
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.7
$ws.Range("S5").Value = 1.42
$ws.Range("T5").Value = 2.47
$ws.Range("W5").Value = 5.3
$ws.Range("X5").Value = 5.7
$ws.Range("AB5").Value = 35
$ws.Range("AC5").Value = 8.5
$ws.Range("AD5").Value = 7.9
$ws.Range("AI5").Value = 65
$ws.Range("AJ5").Value = 28
$ws.Range("AN5").Value = 3
$ws.Range("AO5").Value = 6.5
$ws.Range("AP5").Value = 18.5
$ws.Range("AQ5").Value = 19.5
$ws.Range("AR5").Value = 60
$ws.Range("AT5").Value = 2.45
$ws.Range("AU5").Value = 8.75
$ws.Range("AV5").Value = 100
$ws.Range("AY5").Value = 55
$ws.Range("BA5").Value = 450
$ws.Range("N6").Value = 8
